# Case_0_143: bus voltage magnitude results updated for the 380 kV case
# (external-grid setpoint reduced from 1.05 to 1.02 p.u., propagating through
# the recalculated bus voltages in columns B-F and I-N for rows 2-25).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.066423424578895
$ws.Range("D2").Value = 1.07594709080343
$ws.Range("E2").Value = 1.061311736806519
$ws.Range("F2").Value = 1.083119670135001
$ws.Range("I2").Value = 1.041389914269709
$ws.Range("J2").Value = 1.071373501933423
$ws.Range("K2").Value = 1.078632317410541
$ws.Range("L2").Value = 1.06403614080342
$ws.Range("M2").Value = 1.085786113961695
$ws.Range("N2").Value = 1.072894975194427
$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.068302406642358
$ws.Range("D3").Value = 1.077805316909735
$ws.Range("E3").Value = 1.062944493270285
$ws.Range("F3").Value = 1.08513035300864
$ws.Range("I3").Value = 1.041779487324899
$ws.Range("J3").Value = 1.072904081978044
$ws.Range("K3").Value = 1.080305409465706
$ws.Range("L3").Value = 1.065481432677228
$ws.Range("M3").Value = 1.087612682570042
$ws.Range("N3").Value = 1.074427728838272
$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.069514431094897
$ws.Range("D4").Value = 1.079004231500361
$ws.Range("E4").Value = 1.063997719295265
$ws.Range("F4").Value = 1.086428074374454
$ws.Range("I4").Value = 1.042028666562766
$ws.Range("J4").Value = 1.073890400340859
$ws.Range("K4").Value = 1.081384096031696
$ws.Range("L4").Value = 1.0664128832819
$ws.Range("M4").Value = 1.088790862341889
$ws.Range("N4").Value = 1.075415447886293
$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.070023075806747
$ws.Range("D5").Value = 1.079507440618441
$ws.Range("E5").Value = 1.064439727451967
$ws.Range("F5").Value = 1.086972860973436
$ws.Range("I5").Value = 1.042132732163872
$ws.Range("J5").Value = 1.074304090895797
$ws.Range("K5").Value = 1.081836656133082
$ws.Range("L5").Value = 1.066803582708729
$ws.Range("M5").Value = 1.089285294713678
$ws.Range("N5").Value = 1.075829725929271
$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.070108427856933
$ws.Range("D6").Value = 1.07959188442508
$ws.Range("E6").Value = 1.064513897998716
$ws.Range("F6").Value = 1.087064288159395
$ws.Range("I6").Value = 1.042150164958043
$ws.Range("J6").Value = 1.07437349567744
$ws.Range("K6").Value = 1.081912589543742
$ws.Range("L6").Value = 1.066869131558899
$ws.Range("M6").Value = 1.089368261360324
$ws.Range("N6").Value = 1.075899229273664
$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.069521231101833
$ws.Range("D7").Value = 1.079010958580107
$ws.Range("E7").Value = 1.064003628425592
$ws.Range("F7").Value = 1.086435356858808
$ws.Range("I7").Value = 1.04203005979418
$ws.Range("J7").Value = 1.073895931834895
$ws.Range("K7").Value = 1.081390146751719
$ws.Range("L7").Value = 1.066418107273171
$ws.Range("M7").Value = 1.088797472381285
$ws.Range("N7").Value = 1.075420987235685
$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.067059234801739
$ws.Range("D8").Value = 1.076575819479822
$ws.Range("E8").Value = 1.061864222833341
$ws.Range("F8").Value = 1.083799890686113
$ws.Range("I8").Value = 1.041522175951737
$ws.Range("J8").Value = 1.071891621183509
$ws.Range("K8").Value = 1.07919856824161
$ws.Range("L8").Value = 1.064525369667447
$ws.Range("M8").Value = 1.086404195721812
$ws.Range("N8").Value = 1.073413830233277
$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.062690857370306
$ws.Range("D9").Value = 1.072257267114151
$ws.Range("E9").Value = 1.058068483275229
$ws.Range("F9").Value = 1.0791294710614
$ws.Range("I9").Value = 1.040604776004226
$ws.Range("J9").Value = 1.068327858250897
$ws.Range("K9").Value = 1.075305934923642
$ws.Range("L9").Value = 1.061160712396543
$ws.Range("M9").Value = 1.082157489065545
$ws.Range("N9").Value = 1.069845006348531
$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.059757149637833
$ws.Range("D10").Value = 1.069358517532472
$ws.Range("E10").Value = 1.055519581534449
$ws.Range("F10").Value = 1.075996810603094
$ws.Range("I10").Value = 1.039977763267059
$ws.Range("J10").Value = 1.065929511358373
$ws.Range("K10").Value = 1.072689016133125
$ws.Range("L10").Value = 1.058896858555863
$ws.Range("M10").Value = 1.079305338907699
$ws.Range("N10").Value = 1.067443253528325
$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.058481448297923
$ws.Range("D11").Value = 1.068098380126853
$ws.Range("E11").Value = 1.054411279408068
$ws.Range("F11").Value = 1.074635524158238
$ws.Range("I11").Value = 1.039702531388888
$ws.Range("J11").Value = 1.064885427345862
$ws.Range("K11").Value = 1.0715504308365
$ws.Range("L11").Value = 1.057911444877968
$ws.Range("M11").Value = 1.078065070431861
$ws.Range("N11").Value = 1.066397686796759
$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.058006762105004
$ws.Range("D12").Value = 1.067629540372663
$ws.Range("E12").Value = 1.053998893689268
$ws.Range("F12").Value = 1.07412913173244
$ws.Range("I12").Value = 1.039599731083466
$ws.Range("J12").Value = 1.064496748689955
$ws.Range("K12").Value = 1.071126670004344
$ws.Range("L12").Value = 1.057544625655912
$ws.Range("M12").Value = 1.077603564352484
$ws.Range("N12").Value = 1.066008456172578
$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.058108622040417
$ws.Range("D13").Value = 1.067730143253113
$ws.Range("E13").Value = 1.054087384416632
$ws.Range("F13").Value = 1.074237788977102
$ws.Range("I13").Value = 1.039621807855809
$ws.Range("J13").Value = 1.064580160800222
$ws.Range("K13").Value = 1.071217606488326
$ws.Range("L13").Value = 1.057623345803252
$ws.Range("M13").Value = 1.077702596242288
$ws.Range("N13").Value = 1.06609198673761
$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.058442227744546
$ws.Range("D14").Value = 1.068059641502961
$ws.Range("E14").Value = 1.054377206138163
$ws.Range("F14").Value = 1.074593681050359
$ws.Range("I14").Value = 1.039694045485906
$ws.Range("J14").Value = 1.064853316696954
$ws.Range("K14").Value = 1.071515419888586
$ws.Range("L14").Value = 1.057881139768319
$ws.Range("M14").Value = 1.078026938961708
$ws.Range("N14").Value = 1.066365530547046
$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.05864766198427
$ws.Range("D15").Value = 1.068262553692911
$ws.Range("E15").Value = 1.054555679781667
$ws.Range("F15").Value = 1.074812857753464
$ws.Range("I15").Value = 1.039738478158299
$ws.Range("J15").Value = 1.065021502569689
$ws.Range("K15").Value = 1.071698800642822
$ws.Range("L15").Value = 1.058039869525049
$ws.Range("M15").Value = 1.07822666848994
$ws.Range("N15").Value = 1.066533955263012
$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.059841698149982
$ws.Range("D16").Value = 1.069442042210209
$ws.Range("E16").Value = 1.055593037030381
$ws.Range("F16").Value = 1.076087050916949
$ws.Range("I16").Value = 1.03999595034726
$ws.Range("J16").Value = 1.065998684389704
$ws.Range("K16").Value = 1.072764463763714
$ws.Range("L16").Value = 1.058962147097616
$ws.Range("M16").Value = 1.079387538431257
$ws.Range("N16").Value = 1.067512524793296
$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.060589226066499
$ws.Range("D17").Value = 1.070180559868968
$ws.Range("E17").Value = 1.056242495172352
$ws.Range("F17").Value = 1.076885009398535
$ws.Range("I17").Value = 1.040156452435506
$ws.Range("J17").Value = 1.066610135399346
$ws.Range("K17").Value = 1.0734314536682
$ws.Range("L17").Value = 1.059539275331636
$ws.Range("M17").Value = 1.08011429435798
$ws.Range("N17").Value = 1.068124844133506
$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.061024728412274
$ws.Range("D18").Value = 1.070610847786161
$ws.Range("E18").Value = 1.05662086961877
$ws.Range("F18").Value = 1.077349981380716
$ws.Range("I18").Value = 1.040249711013167
$ws.Range("J18").Value = 1.06696624704952
$ws.Range("K18").Value = 1.073819973998593
$ws.Range("L18").Value = 1.05987540885815
$ws.Range("M18").Value = 1.080537692056767
$ws.Range("N18").Value = 1.068481461503069
$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.061173136329992
$ws.Range("D19").Value = 1.070757484654794
$ws.Range("E19").Value = 1.05674981090974
$ws.Range("F19").Value = 1.077508446919891
$ws.Range("I19").Value = 1.04028144897888
$ws.Range("J19").Value = 1.067087581384855
$ws.Range("K19").Value = 1.073952361239057
$ws.Range("L19").Value = 1.059989938210673
$ws.Range("M19").Value = 1.080681974488069
$ws.Range("N19").Value = 1.068602968147077
$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.060509077156582
$ws.Range("D20").Value = 1.070101373414575
$ws.Range("E20").Value = 1.056172860479207
$ws.Range("F20").Value = 1.076799444150269
$ws.Range("I20").Value = 1.040139269318741
$ws.Range("J20").Value = 1.06654458822359
$ws.Range("K20").Value = 1.07335994627982
$ws.Range("L20").Value = 1.05947740633381
$ws.Range("M20").Value = 1.080036372909225
$ws.Range("N20").Value = 1.068059203873242
$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.058344012399324
$ws.Range("D21").Value = 1.067962633888524
$ws.Range("E21").Value = 1.054291880709893
$ws.Range("F21").Value = 1.074488900614579
$ws.Range("I21").Value = 1.039672789005903
$ws.Range("J21").Value = 1.064772902981811
$ws.Range("K21").Value = 1.071427744618512
$ws.Range("L21").Value = 1.057805247907198
$ws.Range("M21").Value = 1.0779314508079
$ws.Range("N21").Value = 1.066285002635205
$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.056977913028491
$ws.Range("D22").Value = 1.066613464637292
$ws.Range("E22").Value = 1.053105098870996
$ws.Range("F22").Value = 1.073031817966432
$ws.Range("I22").Value = 1.039376211013377
$ws.Range("J22").Value = 1.063653991799135
$ws.Range("K22").Value = 1.07020802424306
$ws.Range("L22").Value = 1.05674929959262
$ws.Range("M22").Value = 1.07660327464406
$ws.Range("N22").Value = 1.065164502470307
$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.057702574982281
$ws.Range("D23").Value = 1.067329115349951
$ws.Range("E23").Value = 1.053734633175139
$ws.Range("F23").Value = 1.073804666126359
$ws.Range("I23").Value = 1.039533745984456
$ws.Range("J23").Value = 1.064247626953047
$ws.Range("K23").Value = 1.070855089969425
$ws.Range("L23").Value = 1.057309519818202
$ws.Range("M23").Value = 1.077307822087182
$ws.Range("N23").Value = 1.065758980654229
$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.060545294599087
$ws.Range("D24").Value = 1.070137155837916
$ws.Range("E24").Value = 1.056204326769577
$ws.Range("F24").Value = 1.076838108835023
$ws.Range("I24").Value = 1.0401470347405
$ws.Range("J24").Value = 1.066574207830838
$ws.Range("K24").Value = 1.073392259007972
$ws.Range("L24").Value = 1.059505363801756
$ws.Range("M24").Value = 1.080071583823031
$ws.Range("N24").Value = 1.068088865543729
$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.063823875431349
$ws.Range("D25").Value = 1.0733771017829
$ws.Range("E25").Value = 1.059052940552648
$ws.Range("F25").Value = 1.080340147163526
$ws.Range("I25").Value = 1.040844638534465
$ws.Range("J25").Value = 1.069253065267704
$ws.Range("K25").Value = 1.076316038092202
$ws.Range("L25").Value = 1.062034141909558
$ws.Range("M25").Value = 1.083258978184648
$ws.Range("N25").Value = 1.070771527265425

Write-Output "Updated 264 cells in vm_pu.xlsx (Case_0_143, 380 kV)"
